{"js": "// Insert a new \"Problem1:\" section (Fisher's Linear Discriminant write-up)\n// at the very top of the document body, before the existing\n// \"Problem 2: Using Low Rank Structure for Corrupted Entries\" paragraph.\n\nconst body = context.document.body;\n\nfunction wrapOoxml(innerXml) {\n  return (\n    '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    \"<pkg:xmlData>\" +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    \"<w:body>\" +\n    innerXml +\n    \"</w:body></w:document>\" +\n    \"</pkg:xmlData></pkg:part></pkg:package>\"\n  );\n}\n\n// All of the new paragraphs to insert, in order, as raw <w:p> OOXML so we\n// can faithfully reproduce plain text paragraphs, truly empty paragraphs\n// (<w:p/>) and the spell-check proofErr spans Word itself generates around\n// \"s_w\", \"sklearn\" and \"lda\".\nconst newParas = [\n  \"<w:p><w:r><w:t>Problem1:</w:t></w:r></w:p>\",\n\n  \"<w:p>\" +\n    '<w:r><w:t xml:space=\"preserve\">Using Fisher\\u2019s Linear Discriminant and the fact that the matrix has full rank we are able to calculate w by the equation </w:t></w:r>' +\n    '<w:proofErr w:type=\"spellStart\"/>' +\n    \"<w:r><w:t>s_w</w:t></w:r>\" +\n    '<w:proofErr w:type=\"spellEnd\"/>' +\n    '<w:r><w:t xml:space=\"preserve\">^-1 * (mean1 \\u2013 mean2) to get </w:t></w:r>' +\n  \"</w:p>\",\n\n  \"<w:p><w:r><w:t>[[ 0.06773847]</w:t></w:r></w:p>\",\n  '<w:p><w:r><w:t xml:space=\"preserve\"> [ 0.12832633]</w:t></w:r></w:p>',\n  '<w:p><w:r><w:t xml:space=\"preserve\"> [-0.17730653]]</w:t></w:r></w:p>',\n\n  \"<w:p/>\",\n  \"<w:p/>\",\n\n  \"<w:p>\" +\n    '<w:r><w:t xml:space=\"preserve\">Using </w:t></w:r>' +\n    '<w:proofErr w:type=\"spellStart\"/>' +\n    \"<w:r><w:t>sklearn</w:t></w:r>\" +\n    '<w:proofErr w:type=\"spellEnd\"/>' +\n    '<w:r><w:t xml:space=\"preserve\"> </w:t></w:r>' +\n    '<w:proofErr w:type=\"spellStart\"/>' +\n    \"<w:r><w:t>lda</w:t></w:r>\" +\n    '<w:proofErr w:type=\"spellEnd\"/>' +\n    '<w:r><w:t xml:space=\"preserve\"> we get w equal to </w:t></w:r>' +\n  \"</w:p>\",\n\n  \"<w:p><w:r><w:t>[[-2.70953878]</w:t></w:r></w:p>\",\n  '<w:p><w:r><w:t xml:space=\"preserve\"> [-5.1330531 ]</w:t></w:r></w:p>',\n  '<w:p><w:r><w:t xml:space=\"preserve\"> [ 7.09226126]]</w:t></w:r></w:p>',\n\n  \"<w:p/>\",\n];\n\n// Create one throw-away anchor paragraph at the very start of the body,\n// then grow the new content after it (in order), and finally remove the\n// throw-away placeholder once every real paragraph is in place.\nlet anchor = body.insertParagraph(\"\", Word.InsertLocation.start);\n\nfor (const paraXml of newParas) {\n  anchor = anchor.insertParagraph(\"\", Word.InsertLocation.after);\n  anchor.insertOoxml(wrapOoxml(paraXml), Word.InsertLocation.replace);\n}\n\n// Remove the placeholder paragraph that was only used as an insertion\n// anchor (it sits before \"Problem1:\" and must not remain in the output).\nconst placeholder = body.paragraphs.getFirst();\nplaceholder.delete();\nawait context.sync();\n\n// Word automatically maintains a \"_GoBack\" bookmark at the location of the\n// most recent edit. Since the edit here happened at the top of the\n// document, re-saving moves \"_GoBack\" from its old spot (end of the\n// corrMat1 paragraph) to the start of the paragraph that now immediately\n// follows the newly-typed text (the \"Problem 2: ...\" paragraph).\ncontext.document.deleteBookmark(\"_GoBack\");\nconst problem2Para = anchor.getNext();\nproblem2Para.getRange(Word.RangeLocation.start).insertBookmark(\"_GoBack\");\n\nawait context.sync();\n", "ps1": "# Insert a new \"Problem1:\" section (Fisher's Linear Discriminant write-up)\n# at the very top of the document body, before the existing\n# \"Problem 2: Using Low Rank Structure for Corrupted Entries\" paragraph.\n\n$d = $word.ActiveDocument\n\n# Build the literal characters Word itself would use: U+2019 (right single\n# quotation mark) and U+2013 (en dash).\n$rsquo = [char]0x2019\n$ndash = [char]0x2013\n\n$ns = 'xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"'\n\n$body = \"\"\n$body += \"<w:p><w:r><w:t>Problem1:</w:t></w:r></w:p>\"\n$body += \"<w:p>\"\n$body += \"<w:r><w:t xml:space=`\"preserve`\">Using Fisher$($rsquo)s Linear Discriminant and the fact that the matrix has full rank we are able to calculate w by the equation </w:t></w:r>\"\n$body += \"<w:proofErr w:type=`\"spellStart`\"/>\"\n$body += \"<w:r><w:t>s_w</w:t></w:r>\"\n$body += \"<w:proofErr w:type=`\"spellEnd`\"/>\"\n$body += \"<w:r><w:t xml:space=`\"preserve`\">^-1 * (mean1 $($ndash) mean2) to get </w:t></w:r>\"\n$body += \"</w:p>\"\n$body += \"<w:p><w:r><w:t>[[ 0.06773847]</w:t></w:r></w:p>\"\n$body += \"<w:p><w:r><w:t xml:space=`\"preserve`\"> [ 0.12832633]</w:t></w:r></w:p>\"\n$body += \"<w:p><w:r><w:t xml:space=`\"preserve`\"> [-0.17730653]]</w:t></w:r></w:p>\"\n$body += \"<w:p/>\"\n$body += \"<w:p/>\"\n$body += \"<w:p>\"\n$body += \"<w:r><w:t xml:space=`\"preserve`\">Using </w:t></w:r>\"\n$body += \"<w:proofErr w:type=`\"spellStart`\"/>\"\n$body += \"<w:r><w:t>sklearn</w:t></w:r>\"\n$body += \"<w:proofErr w:type=`\"spellEnd`\"/>\"\n$body += \"<w:r><w:t xml:space=`\"preserve`\"> </w:t></w:r>\"\n$body += \"<w:proofErr w:type=`\"spellStart`\"/>\"\n$body += \"<w:r><w:t>lda</w:t></w:r>\"\n$body += \"<w:proofErr w:type=`\"spellEnd`\"/>\"\n$body += \"<w:r><w:t xml:space=`\"preserve`\"> we get w equal to </w:t></w:r>\"\n$body += \"</w:p>\"\n$body += \"<w:p><w:r><w:t>[[-2.70953878]</w:t></w:r></w:p>\"\n$body += \"<w:p><w:r><w:t xml:space=`\"preserve`\"> [-5.1330531 ]</w:t></w:r></w:p>\"\n$body += \"<w:p><w:r><w:t xml:space=`\"preserve`\"> [ 7.09226126]]</w:t></w:r></w:p>\"\n$body += \"<w:p/>\"\n\n$xml = \"<?xml version=`\"1.0`\" encoding=`\"UTF-8`\" standalone=`\"yes`\"?><w:body $ns>$body</w:body>\"\n\n# Insert all of the new paragraphs as a single block at the very start of\n# the document (collapsed range at position 0).\n$r = $d.Range(0, 0)\n$r.InsertXML($xml)\n\nWrite-Output \"done\"\n"}
